$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be stored as text so numeric-looking
# strings (e.g. "490.10", "0.333") keep their exact literal formatting
# instead of being auto-coerced into floating point numbers by the
# COM value setter. ClearFormats() afterwards drops the temporary
# "@" number format again so cell formatting is left untouched.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "56.796.06"
$ws.Range("E2").Value = "  +3.62%  "

$ws.Range("D3").Value = "2.469.57"
$ws.Range("E3").Value = "  +1.74%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "490.10"
$ws.Range("E5").Value = "  +2.59%  "

$ws.Range("D6").Value = "151.85"
$ws.Range("E6").Value = "  +9.30%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("E8").Value = "  +2.28%  "

$ws.Range("D9").Value = "2.475.49"
$ws.Range("E9").Value = "  +0.99%  "

$ws.Range("D10").Value = "0.0994"
$ws.Range("E10").Value = "  +3.80%  "

$ws.Range("E11").Value = "  +4.26%  "

$ws.Range("D12").Value = "0.333"
$ws.Range("E12").Value = "  +3.12%  "

$ws.Range("E13").Value = "  +1.43%  "

$ws.Range("D14").Value = "2.902.27"
$ws.Range("E14").Value = "  +1.63%  "

$ws.Range("D15").Value = "57.089.01"
$ws.Range("E15").Value = "  +3.87%  "

$ws.Range("D16").Value = "20.99"
$ws.Range("E16").Value = "  +2.85%  "

$ws.Range("D17").Value = "0.0000137"
$ws.Range("E17").Value = "  +2.33%  "

$ws.Range("D18").Value = "2.462.99"
$ws.Range("E18").Value = "  +0.53%  "

$ws.Range("D19").Value = "4.55"
$ws.Range("E19").Value = "  +5.04%  "

$ws.Range("D20").Value = "10.14"
$ws.Range("E20").Value = "  +2.98%  "

$ws.Range("D21").Value = "320.45"
$ws.Range("E21").Value = "  +2.37%  "

$ws.Range("D22").Value = "0.997"
$ws.Range("E22").Value = "  +0.45%  "

$ws.Range("E23").Value = "  +4.07%  "

$ws.Range("D24").Value = "58.06"
$ws.Range("E24").Value = "  +1.74%  "

$ws.Range("E25").Value = "  +1.28%  "

$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("E27").Value = "  +0.60%  "

$ws.Range("D28").Value = "2.586.01"
$ws.Range("E28").Value = "  +1.75%  "

$ws.Range("D29").Value = "7.58"
$ws.Range("E29").Value = "  +3.34%  "

$ws.Range("D30").Value = "0.0₃0805"
$ws.Range("E30").Value = "  +5.20%  "

$ws.Range("E31").Value = "  +0.23%  "

$ws.Range("D32").Value = "150.72"
$ws.Range("E32").Value = "  +1.47%  "

$ws.Range("D33").Value = "18.26"
$ws.Range("E33").Value = "  +2.29%  "

$ws.Range("D34").Value = "1.52"
$ws.Range("E34").Value = "  +3.44%  "

$ws.Range("E35").Value = "  +1.34%  "

$ws.Range("E36").Value = "  +3.00%  "

$ws.Range("D37").Value = "0.887"
$ws.Range("E37").Value = "  +6.07%  "

$ws.Range("D38").Value = "3.76"
$ws.Range("E38").Value = "  +5.35%  "

$ws.Range("D39").Value = "34.11"
$ws.Range("E39").Value = "  +1.93%  "

$ws.Range("E40").Value = "  +8.17%  "

$ws.Range("E41").Value = "  +3.05%  "

$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "0.995"
$ws.Range("E42").Value = "  +0.33%  "

$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").Value = "0.0556"
$ws.Range("E43").Value = "  +2.81%  "

$ws.Range("E44").Value = "  +1.25%  "

$ws.Range("D45").Value = "0.0947"
$ws.Range("E45").Value = "  +6.23%  "

$ws.Range("D46").Value = "4.80"
$ws.Range("E46").Value = "  +4.01%  "

$ws.Range("D47").Value = "262.51"
$ws.Range("E47").Value = "  +4.61%  "

$ws.Range("E48").Value = "  +0.73%  "

$ws.Range("E49").Value = "  +3.13%  "

$ws.Range("D50").Value = "17.69"
$ws.Range("E50").Value = "  +3.21%  "

$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "1.73"
$ws.Range("E51").Value = "  +25.78%  "

$priceRange.ClearFormats()
